$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the Espárragos price rows: dates and associated
# price/volume/unit figures are rotated among the existing data rows.
$ws.Range("D2").Value = 44511
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 1300
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = 1350
$ws.Range("P2").Value = 1350
$ws.Range("D3").Value = 44839
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1700
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = 1760
$ws.Range("P3").Value = 1760
$ws.Range("D5").Value = 44468
$ws.Range("H5").Value = "Verde"
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1920
$ws.Range("P5").Value = 1920
$ws.Range("D6").Value = 44477
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("K6").Value = 1400
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 1460
$ws.Range("P6").Value = 1460
$ws.Range("D7").Value = 44545
$ws.Range("J7").Value = 550
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = 1755
$ws.Range("P7").Value = 1755
$ws.Range("D8").Value = 44496
$ws.Range("J8").Value = 550
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 1773
$ws.Range("N8").Value = "$/paquete"
$ws.Range("P8").Value = 1773
$ws.Range("D9").Value = 44519
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 1300
$ws.Range("M9").Value = 1240
$ws.Range("P9").Value = 1240
$ws.Range("D10").Value = 44510
$ws.Range("D12").Value = 44489
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1450
$ws.Range("N12").Value = "$/kilo"
$ws.Range("P12").Value = 1450
